$wb = $excel.ActiveWorkbook

# The "NegativeLogins" sheet has validation-message cells that used to contain
# the specific messages "Username cannot be empty" / "Password cannot be empty".
# The web elements under test changed, so these now just show "Required".
$ws2 = $wb.Worksheets.Item("NegativeLogins")
$ws2.Range("C5").Value = "Required"
$ws2.Range("C6").Value = "Required"
$ws2.Range("C7").Value = "Required"

# The active/selected sheet moved from Employee to NegativeLogins, with a new
# selected cell on each sheet.
$ws1 = $wb.Worksheets.Item("Employee")
$ws1.Range("D11").Select()

$ws2.Activate()
$ws2.Range("C11").Select()
